$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update corrected system addresses (BIOS/BDOS/CCP target address column)
$ws.Range("D8").Value = '$CC00 .. $CFFF'
$ws.Range("D9").Value = '$D000 .. $D7FF'
$ws.Range("D10").Value = '$D800 .. $DFFF'

# Update the selected/active cell to match the saved view state
$ws.Range("E20").Select()
